$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "YYYY [YRYYYY]" text headers in E1:BL1 with plain numeric years
# (1960-2019), left-aligned - the "2020 [YR2020]" header in BM1 is left as-is.
$col = 5
for ($year = 1960; $year -le 2019; $year++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $year
    $cell.HorizontalAlignment = -4131
    $col++
}

# Update the visible selection/scroll position to match the edited range.
$ws.Range("E1:BL1").Select() | Out-Null
